$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): swap column order of category labels
$ws.Range("B1").Value = "bedrooms_2"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "living_rooms_1"
$ws.Range("E1").Value = "kitchens_2"

# Row 3: move the "1" from D3 to E3
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1

# Row 5: move the "1" from C5 to B5
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0

# Row 6: move the "1" from B6 to C6
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1

# Row 7: move the "1" from E7 to D7
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
